$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 118  # H12: 120.75 -> 118
$ws.Cells.Item(12, 9).Value = 119.5  # I12: 122.63636 -> 119.5
$ws.Cells.Item(12, 11).Value = 119.5  # K12: 122.63636 -> 119.5
$ws.Cells.Item(12, 13).Value = 50.5  # M12: 47.36364 -> 50.5

$ws.Cells.Item(17, 8).Value = 402.32144  # H17: 382.55554 -> 402.32144
$ws.Cells.Item(17, 10).Value = 402.32144  # J17: 382.55554 -> 402.32144
$ws.Cells.Item(17, 12).Value = 1206.96432  # L17: 1147.66662 -> 1206.96432
$ws.Cells.Item(17, 14).Value = -1542.96432  # N17: -1483.66662 -> -1542.96432

$ws.Cells.Item(32, 8).Value = 18225.385  # H32: 18379.385 -> 18225.385
$ws.Cells.Item(32, 9).Value = 17999.5  # I32: 19000.5 -> 17999.5
$ws.Cells.Item(32, 11).Value = 17999.5  # K32: 19000.5 -> 17999.5
$ws.Cells.Item(32, 13).Value = -17673.5  # M32: -18674.5 -> -17673.5

$ws.Cells.Item(55, 8).Value = 258.54544  # H55: 264.875 -> 258.54544
$ws.Cells.Item(55, 9).Value = 189.4  # I55: 279.4 -> 189.4
$ws.Cells.Item(55, 10).Value = 316.16666  # J55: 258.27274 -> 316.16666
$ws.Cells.Item(55, 11).Value = 189.4  # K55: 279.4 -> 189.4
$ws.Cells.Item(55, 12).Value = 316.16666  # L55: 258.27274 -> 316.16666
$ws.Cells.Item(55, 13).Value = 24.59999999999999  # M55: -65.39999999999998 -> 24.59999999999999
$ws.Cells.Item(55, 14).Value = -744.16666  # N55: -686.27274 -> -744.16666

$ws.Cells.Item(70, 8).Value = 1498.5  # H70: 1499.5 -> 1498.5
$ws.Cells.Item(70, 9).Value = 1497  # I70: 0 -> 1497
$ws.Cells.Item(70, 10).Value = 1499  # J70: 1499.5 -> 1499
$ws.Cells.Item(70, 11).Value = 4491  # K70: 0 -> 4491
$ws.Cells.Item(70, 12).Value = 4497  # L70: 4498.5 -> 4497
$ws.Cells.Item(70, 13).Value = -4221  # M70: None -> -4221
$ws.Cells.Item(70, 14).Value = -5037  # N70: -5038.5 -> -5037

$ws.Cells.Item(73, 8).Value = 1498.5  # H73: 1499.5 -> 1498.5
$ws.Cells.Item(73, 9).Value = 1497  # I73: 0 -> 1497
$ws.Cells.Item(73, 10).Value = 1499  # J73: 1499.5 -> 1499
$ws.Cells.Item(73, 11).Value = 4491  # K73: 0 -> 4491
$ws.Cells.Item(73, 12).Value = 4497  # L73: 4498.5 -> 4497
$ws.Cells.Item(73, 13).Value = -3555  # M73: None -> -3555
$ws.Cells.Item(73, 14).Value = -6369  # N73: -6370.5 -> -6369

$ws.Cells.Item(111, 8).Value = 3621.1667  # H111: 3671.1667 -> 3621.1667
$ws.Cells.Item(111, 10).Value = 3449  # J111: 3599 -> 3449
$ws.Cells.Item(111, 12).Value = 10347  # L111: 10797 -> 10347
$ws.Cells.Item(111, 14).Value = -16481  # N111: -16931 -> -16481

$ws.Cells.Item(132, 8).Value = 6165.8613  # H132: 6772.75 -> 6165.8613
$ws.Cells.Item(132, 9).Value = 3159.9285  # I132: 3468.125 -> 3159.9285
$ws.Cells.Item(132, 11).Value = 9479.7855  # K132: 10404.375 -> 9479.7855
$ws.Cells.Item(132, 13).Value = -6949.7855  # M132: -7874.375 -> -6949.7855

$ws.Cells.Item(138, 8).Value = 4318.294  # H138: 3661.524 -> 4318.294
$ws.Cells.Item(138, 9).Value = 1205.1666  # I138: 1013.7857 -> 1205.1666
$ws.Cells.Item(138, 11).Value = 3615.4998  # K138: 3041.3571 -> 3615.4998
$ws.Cells.Item(138, 13).Value = 1524.5002  # M138: 2098.6429 -> 1524.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5301.6787  # H2: 5659.654 -> 5301.6787
$ws.Cells.Item(2, 9).Value = 3993.238  # I2: 4345.3687 -> 3993.238
$ws.Cells.Item(2, 11).Value = 3993.238  # K2: 4345.3687 -> 3993.238
$ws.Cells.Item(2, 13).Value = -3880.238  # M2: -4232.3687 -> -3880.238

$ws.Cells.Item(4, 8).Value = 450  # H4: 391.66666 -> 450
$ws.Cells.Item(4, 9).Value = 496.75  # I4: 417.4 -> 496.75
$ws.Cells.Item(4, 11).Value = 496.75  # K4: 417.4 -> 496.75
$ws.Cells.Item(4, 13).Value = -380.75  # M4: -301.4 -> -380.75

$ws.Cells.Item(5, 8).Value = 0  # H5: 222 -> 0
$ws.Cells.Item(5, 9).Value = 0  # I5: 222 -> 0
$ws.Cells.Item(5, 11).Value = 0  # K5: 222 -> 0
$ws.Cells.Item(5, 13).ClearContents()  # was M5=-110, now removed

$ws.Cells.Item(27, 8).Value = 6008  # H27: 6004 -> 6008
$ws.Cells.Item(27, 10).Value = 6008  # J27: 6004 -> 6008
$ws.Cells.Item(27, 12).Value = 6008  # L27: 6004 -> 6008
$ws.Cells.Item(27, 14).Value = -6376  # N27: -6372 -> -6376

$ws.Cells.Item(33, 8).Value = 4097  # H33: 4055.5557 -> 4097
$ws.Cells.Item(33, 9).Value = 3996.6667  # I33: 3857.1428 -> 3996.6667
$ws.Cells.Item(33, 10).Value = 5000  # J33: 4750 -> 5000
$ws.Cells.Item(33, 11).Value = 3996.6667  # K33: 3857.1428 -> 3996.6667
$ws.Cells.Item(33, 12).Value = 5000  # L33: 4750 -> 5000
$ws.Cells.Item(33, 13).Value = -3667.6667  # M33: -3528.1428 -> -3667.6667
$ws.Cells.Item(33, 14).Value = -5658  # N33: -5408 -> -5658

$ws.Cells.Item(45, 8).Value = 1410.6666  # H45: 1431.0769 -> 1410.6666
$ws.Cells.Item(45, 9).Value = 1287.25  # I45: 1304.9565 -> 1287.25
$ws.Cells.Item(45, 11).Value = 1287.25  # K45: 1304.9565 -> 1287.25
$ws.Cells.Item(45, 13).Value = -910.25  # M45: -927.9565 -> -910.25

$ws.Cells.Item(63, 8).Value = 2428.4783  # H63: 2657.75 -> 2428.4783
$ws.Cells.Item(63, 9).Value = 1138  # I63: 1275.6 -> 1138
$ws.Cells.Item(63, 10).Value = 3836.2727  # J63: 4039.9 -> 3836.2727
$ws.Cells.Item(63, 11).Value = 1138  # K63: 1275.6 -> 1138
$ws.Cells.Item(63, 12).Value = 3836.2727  # L63: 4039.9 -> 3836.2727
$ws.Cells.Item(63, 13).Value = -452  # M63: -589.5999999999999 -> -452
$ws.Cells.Item(63, 14).Value = -5208.2727  # N63: -5411.9 -> -5208.2727

$ws.Cells.Item(66, 8).Value = 2428.4783  # H66: 2657.75 -> 2428.4783
$ws.Cells.Item(66, 9).Value = 1138  # I66: 1275.6 -> 1138
$ws.Cells.Item(66, 10).Value = 3836.2727  # J66: 4039.9 -> 3836.2727
$ws.Cells.Item(66, 11).Value = 5690  # K66: 6378 -> 5690
$ws.Cells.Item(66, 12).Value = 19181.3635  # L66: 20199.5 -> 19181.3635
$ws.Cells.Item(66, 13).Value = -2258  # M66: -2946 -> -2258
$ws.Cells.Item(66, 14).Value = -26045.3635  # N66: -27063.5 -> -26045.3635

$ws.Cells.Item(88, 8).Value = 7134.6875  # H88: 8528.308000000001 -> 7134.6875
$ws.Cells.Item(88, 9).Value = 1535  # I88: 1921.5 -> 1535
$ws.Cells.Item(88, 10).Value = 9001.25  # J88: 9729.546 -> 9001.25
$ws.Cells.Item(88, 11).Value = 1535  # K88: 1921.5 -> 1535
$ws.Cells.Item(88, 12).Value = 9001.25  # L88: 9729.546 -> 9001.25
$ws.Cells.Item(88, 13).Value = -1129  # M88: -1515.5 -> -1129
$ws.Cells.Item(88, 14).Value = -9813.25  # N88: -10541.546 -> -9813.25

$ws.Cells.Item(91, 8).Value = 7134.6875  # H91: 8528.308000000001 -> 7134.6875
$ws.Cells.Item(91, 9).Value = 1535  # I91: 1921.5 -> 1535
$ws.Cells.Item(91, 10).Value = 9001.25  # J91: 9729.546 -> 9001.25
$ws.Cells.Item(91, 11).Value = 1535  # K91: 1921.5 -> 1535
$ws.Cells.Item(91, 12).Value = 9001.25  # L91: 9729.546 -> 9001.25
$ws.Cells.Item(91, 13).Value = -131  # M91: -517.5 -> -131
$ws.Cells.Item(91, 14).Value = -11809.25  # N91: -12537.546 -> -11809.25

$ws.Cells.Item(116, 8).Value = 5301.6787  # H116: 5659.654 -> 5301.6787
$ws.Cells.Item(116, 9).Value = 3993.238  # I116: 4345.3687 -> 3993.238
$ws.Cells.Item(116, 11).Value = 3993.238  # K116: 4345.3687 -> 3993.238
$ws.Cells.Item(116, 13).Value = -1699.238  # M116: -2051.3687 -> -1699.238

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5301.6787  # H3: 5659.654 -> 5301.6787
$ws.Cells.Item(3, 9).Value = 3993.238  # I3: 4345.3687 -> 3993.238
$ws.Cells.Item(3, 11).Value = 3993.238  # K3: 4345.3687 -> 3993.238
$ws.Cells.Item(3, 13).Value = -3879.238  # M3: -4231.3687 -> -3879.238

$ws.Cells.Item(4, 8).Value = 0  # H4: 222 -> 0
$ws.Cells.Item(4, 9).Value = 0  # I4: 222 -> 0
$ws.Cells.Item(4, 11).Value = 0  # K4: 222 -> 0
$ws.Cells.Item(4, 13).ClearContents()  # was M4=-107, now removed

$ws.Cells.Item(99, 8).Value = 6978.646  # H99: 7710.5815 -> 6978.646
$ws.Cells.Item(99, 10).Value = 4810.1  # J99: 8936.200000000001 -> 4810.1
$ws.Cells.Item(99, 12).Value = 4810.1  # L99: 8936.200000000001 -> 4810.1
$ws.Cells.Item(99, 14).Value = -7806.1  # N99: -11932.2 -> -7806.1

$ws.Cells.Item(105, 8).Value = 1609.8695  # H105: 1735.3334 -> 1609.8695
$ws.Cells.Item(105, 9).Value = 1539.35  # I105: 1672.3334 -> 1539.35
$ws.Cells.Item(105, 10).Value = 2080  # J105: 2113.3333 -> 2080
$ws.Cells.Item(105, 11).Value = 1539.35  # K105: 1672.3334 -> 1539.35
$ws.Cells.Item(105, 12).Value = 2080  # L105: 2113.3333 -> 2080
$ws.Cells.Item(105, 13).Value = 207.6500000000001  # M105: 74.66660000000002 -> 207.6500000000001
$ws.Cells.Item(105, 14).Value = -5574  # N105: -5607.3333 -> -5574

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 184.2  # H7: 175.8125 -> 184.2
$ws.Cells.Item(7, 10).Value = 0  # J7: 50 -> 0
$ws.Cells.Item(7, 12).Value = 0  # L7: 50 -> 0
$ws.Cells.Item(7, 14).ClearContents()  # was N7=-276, now removed

$ws.Cells.Item(16, 8).Value = 4465.8184  # H16: 4761.067 -> 4465.8184
$ws.Cells.Item(16, 9).Value = 4319.727  # I16: 4762.8423 -> 4319.727
$ws.Cells.Item(16, 11).Value = 4319.727  # K16: 4762.8423 -> 4319.727
$ws.Cells.Item(16, 13).Value = -4032.727  # M16: -4475.8423 -> -4032.727

$ws.Cells.Item(22, 8).Value = 3269  # H22: 2066.4 -> 3269
$ws.Cells.Item(22, 9).Value = 1581.6  # I22: 580.55554 -> 1581.6
$ws.Cells.Item(22, 10).Value = 4956.4  # J22: 4295.1665 -> 4956.4
$ws.Cells.Item(22, 11).Value = 1581.6  # K22: 580.55554 -> 1581.6
$ws.Cells.Item(22, 12).Value = 4956.4  # L22: 4295.1665 -> 4956.4
$ws.Cells.Item(22, 13).Value = -1231.6  # M22: -230.55554 -> -1231.6
$ws.Cells.Item(22, 14).Value = -5656.4  # N22: -4995.1665 -> -5656.4

$ws.Cells.Item(99, 8).Value = 2201.6  # H99: 6801.8335 -> 2201.6
$ws.Cells.Item(99, 9).Value = 2201.6  # I99: 2167.8333 -> 2201.6
$ws.Cells.Item(99, 10).Value = 0  # J99: 11435.833 -> 0
$ws.Cells.Item(99, 11).Value = 2201.6  # K99: 2167.8333 -> 2201.6
$ws.Cells.Item(99, 12).Value = 0  # L99: 11435.833 -> 0
$ws.Cells.Item(99, 13).ClearContents()  # was M99=-669.8332999999998, now removed
$ws.Cells.Item(99, 14).Value = -703.5999999999999  # N99: -14431.833 -> -703.5999999999999

$ws.Cells.Item(107, 8).Value = 994.875  # H107: 1021.5217 -> 994.875
$ws.Cells.Item(107, 9).Value = 729.7368  # I107: 749.05554 -> 729.7368
$ws.Cells.Item(107, 11).Value = 729.7368  # K107: 749.05554 -> 729.7368
$ws.Cells.Item(107, 13).Value = 1190.2632  # M107: 1170.94446 -> 1190.2632

$ws.Cells.Item(113, 8).Value = 4465.8184  # H113: 4761.067 -> 4465.8184
$ws.Cells.Item(113, 9).Value = 4319.727  # I113: 4762.8423 -> 4319.727
$ws.Cells.Item(113, 11).Value = 4319.727  # K113: 4762.8423 -> 4319.727
$ws.Cells.Item(113, 13).Value = -2149.727  # M113: -2592.8423 -> -2149.727

$ws.Cells.Item(126, 8).Value = 2201.6  # H126: 6801.8335 -> 2201.6
$ws.Cells.Item(126, 9).Value = 2201.6  # I126: 2167.8333 -> 2201.6
$ws.Cells.Item(126, 10).Value = 0  # J126: 11435.833 -> 0
$ws.Cells.Item(126, 11).Value = 6604.799999999999  # K126: 6503.499899999999 -> 6604.799999999999
$ws.Cells.Item(126, 12).Value = 0  # L126: 34307.499 -> 0
$ws.Cells.Item(126, 13).ClearContents()  # was M126=-4033.499899999999, now removed
$ws.Cells.Item(126, 14).Value = -4134.799999999999  # N126: -39247.499 -> -4134.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 18092.688  # H3: 10945.5 -> 18092.688
$ws.Cells.Item(3, 9).Value = 5385.375  # I3: 6079.2856 -> 5385.375
$ws.Cells.Item(3, 10).Value = 30800  # J3: 22300 -> 30800
$ws.Cells.Item(3, 11).Value = 16156.125  # K3: 18237.8568 -> 16156.125
$ws.Cells.Item(3, 12).Value = 92400  # L3: 66900 -> 92400
$ws.Cells.Item(3, 13).Value = -16044.125  # M3: -18125.8568 -> -16044.125
$ws.Cells.Item(3, 14).Value = -92624  # N3: -67124 -> -92624

$ws.Cells.Item(12, 8).Value = 1250094.8  # H12: 1428678.4 -> 1250094.8
$ws.Cells.Item(12, 10).Value = 160.14285  # J12: 220.4 -> 160.14285
$ws.Cells.Item(12, 12).Value = 480.42855  # L12: 661.2 -> 480.42855
$ws.Cells.Item(12, 14).Value = -826.4285500000001  # N12: -1007.2 -> -826.4285500000001

$ws.Cells.Item(16, 8).Value = 33.333332  # H16: 37.5 -> 33.333332
$ws.Cells.Item(16, 9).Value = 33.333332  # I16: 37.5 -> 33.333332
$ws.Cells.Item(16, 11).Value = 99.999996  # K16: 112.5 -> 99.999996
$ws.Cells.Item(16, 13).Value = 73.000004  # M16: 60.5 -> 73.000004

$ws.Cells.Item(46, 8).Value = 2512  # H46: 4376.5 -> 2512
$ws.Cells.Item(46, 10).Value = 2512  # J46: 4376.5 -> 2512
$ws.Cells.Item(46, 12).Value = 7536  # L46: 13129.5 -> 7536
$ws.Cells.Item(46, 14).Value = -7718  # N46: -13311.5 -> -7718

$ws.Cells.Item(132, 8).Value = 2528.077  # H132: 2413.6667 -> 2528.077
$ws.Cells.Item(132, 9).Value = 1111  # I132: 1154.75 -> 1111
$ws.Cells.Item(132, 10).Value = 4795.4  # J132: 3852.4285 -> 4795.4
$ws.Cells.Item(132, 11).Value = 9999  # K132: 10392.75 -> 9999
$ws.Cells.Item(132, 12).Value = 43158.6  # L132: 34671.8565 -> 43158.6
$ws.Cells.Item(132, 13).Value = -7469  # M132: -7862.75 -> -7469
$ws.Cells.Item(132, 14).Value = -48218.6  # N132: -39731.8565 -> -48218.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 72  # H2: 83 -> 72
$ws.Cells.Item(2, 9).Value = 78.333336  # I2: 95 -> 78.333336
$ws.Cells.Item(2, 10).Value = 53  # J2: 59 -> 53
$ws.Cells.Item(2, 11).Value = 78.333336  # K2: 95 -> 78.333336
$ws.Cells.Item(2, 12).Value = 53  # L2: 59 -> 53
$ws.Cells.Item(2, 13).Value = 34.666664  # M2: 18 -> 34.666664
$ws.Cells.Item(2, 14).Value = -279  # N2: -285 -> -279

$ws.Cells.Item(102, 8).Value = 6155.5713  # H102: 6343.085 -> 6155.5713
$ws.Cells.Item(102, 9).Value = 5193.8335  # I102: 5439.893 -> 5193.8335
$ws.Cells.Item(102, 11).Value = 5193.8335  # K102: 5439.893 -> 5193.8335
$ws.Cells.Item(102, 13).Value = -3571.8335  # M102: -3817.893 -> -3571.8335

$ws.Cells.Item(113, 8).Value = 8788.786  # H113: 9310.846 -> 8788.786
$ws.Cells.Item(113, 9).Value = 2878.8  # I113: 3098 -> 2878.8
$ws.Cells.Item(113, 11).Value = 2878.8  # K113: 3098 -> 2878.8
$ws.Cells.Item(113, 13).Value = -708.8000000000002  # M113: -928 -> -708.8000000000002

$ws.Cells.Item(122, 8).Value = 59437.3  # H122: 62489.26 -> 59437.3
$ws.Cells.Item(122, 9).Value = 135579.38  # I122: 154740.72 -> 135579.38
$ws.Cells.Item(122, 11).Value = 406738.14  # K122: 464222.16 -> 406738.14
$ws.Cells.Item(122, 13).Value = -404288.14  # M122: -461772.16 -> -404288.14

$ws.Cells.Item(126, 8).Value = 8386.799999999999  # H126: 7695.136 -> 8386.799999999999
$ws.Cells.Item(126, 9).Value = 3478.5833  # I126: 3092.8572 -> 3478.5833
$ws.Cells.Item(126, 11).Value = 10435.7499  # K126: 9278.571599999999 -> 10435.7499
$ws.Cells.Item(126, 13).Value = -7965.749899999999  # M126: -6808.571599999999 -> -7965.749899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5691.4546  # H7: 5824.524 -> 5691.4546
$ws.Cells.Item(7, 9).Value = 4557.5  # I7: 4794.7144 -> 4557.5
$ws.Cells.Item(7, 11).Value = 4557.5  # K7: 4794.7144 -> 4557.5
$ws.Cells.Item(7, 13).Value = -4445.5  # M7: -4682.7144 -> -4445.5

$ws.Cells.Item(40, 8).Value = 4328.3125  # H40: 4772.0713 -> 4328.3125
$ws.Cells.Item(40, 9).Value = 2604.8572  # I40: 3158 -> 2604.8572
$ws.Cells.Item(40, 11).Value = 2604.8572  # K40: 3158 -> 2604.8572
$ws.Cells.Item(40, 13).Value = -2468.8572  # M40: -3022 -> -2468.8572

$ws.Cells.Item(88, 8).Value = 15674.333  # H88: 20062.5 -> 15674.333
$ws.Cells.Item(88, 9).Value = 11917.75  # I88: 22000 -> 11917.75
$ws.Cells.Item(88, 10).Value = 23187.5  # J88: 19093.75 -> 23187.5
$ws.Cells.Item(88, 11).Value = 11917.75  # K88: 22000 -> 11917.75
$ws.Cells.Item(88, 12).Value = 23187.5  # L88: 19093.75 -> 23187.5
$ws.Cells.Item(88, 13).Value = -11489.75  # M88: -21572 -> -11489.75
$ws.Cells.Item(88, 14).Value = -24043.5  # N88: -19949.75 -> -24043.5

$ws.Cells.Item(91, 8).Value = 15674.333  # H91: 20062.5 -> 15674.333
$ws.Cells.Item(91, 9).Value = 11917.75  # I91: 22000 -> 11917.75
$ws.Cells.Item(91, 10).Value = 23187.5  # J91: 19093.75 -> 23187.5
$ws.Cells.Item(91, 11).Value = 11917.75  # K91: 22000 -> 11917.75
$ws.Cells.Item(91, 12).Value = 23187.5  # L91: 19093.75 -> 23187.5
$ws.Cells.Item(91, 13).Value = -10435.75  # M91: -20518 -> -10435.75
$ws.Cells.Item(91, 14).Value = -26151.5  # N91: -22057.75 -> -26151.5

$ws.Cells.Item(126, 8).Value = 5691.4546  # H126: 5824.524 -> 5691.4546
$ws.Cells.Item(126, 9).Value = 4557.5  # I126: 4794.7144 -> 4557.5
$ws.Cells.Item(126, 11).Value = 13672.5  # K126: 14384.1432 -> 13672.5
$ws.Cells.Item(126, 13).Value = -11202.5  # M126: -11914.1432 -> -11202.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3633.35  # H122: 3706.6843 -> 3633.35
$ws.Cells.Item(122, 10).Value = 2241.5  # J122: 2243 -> 2241.5
$ws.Cells.Item(122, 12).Value = 6724.5  # L122: 6729 -> 6724.5
$ws.Cells.Item(122, 14).Value = -11624.5  # N122: -11629 -> -11624.5
